$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'256.70"
$ws.Range("E2").Value = "'0.78%"
$ws.Range("D3").Value = "'26.89"
$ws.Range("E3").Value = "'-4.80%"
$ws.Range("D4").Value = "'4.715"
$ws.Range("E4").Value = "'-10.06%"
$ws.Range("D5").Value = "'0.05933"
$ws.Range("E5").Value = "'0.78%"
$ws.Range("D6").Value = "'6.661"
$ws.Range("E6").Value = "'-1.09%"
$ws.Range("D7").Value = "'0.8673"
$ws.Range("E7").Value = "'0.18%"
$ws.Range("D8").Value = "'0.9399"
$ws.Range("E8").Value = "'-5.02%"
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").Value = "'0.1399"
$ws.Range("E9").Value = "'-0.76%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.03730"
$ws.Range("E10").Value = "'7.89%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.07118"
$ws.Range("E11").Value = "'-0.74%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.03160"
$ws.Range("E12").Value = "'-0.85%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.09246"
$ws.Range("E13").Value = "'0.22%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001534"
$ws.Range("E14").Value = "'-0.70%"
$ws.Range("B15").Value = "One"
$ws.Range("C15").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D15").Value = "'0.0006054"
$ws.Range("E15").Value = "'-0.79%"
$ws.Range("D16").Value = "'0.006000"
$ws.Range("E16").Value = "'2.77%"
$ws.Range("D17").Value = "'3.492"
$ws.Range("E17").Value = "'-0.20%"
$ws.Range("D18").Value = "'3.203"
$ws.Range("E18").Value = "'-0.62%"
$ws.Range("E19").Value = "'1.68%"
$ws.Range("D20").Value = "'0.3126"
$ws.Range("E20").Value = "'-1.66%"
$ws.Range("E21").Value = "'0.34%"
$ws.Range("D22").Value = "'3.805"
$ws.Range("E22").Value = "'7.67%"
$ws.Range("D23").Value = "'0.04215"
$ws.Range("E23").Value = "'1.62%"
$ws.Range("D25").Value = "'0.001224"
$ws.Range("E25").Value = "'-0.17%"
$ws.Range("D26").Value = "'0.004286"
$ws.Range("E26").Value = "'-10.61%"
$ws.Range("E27").Value = "'-0.13%"
$ws.Range("D28").Value = "'0.0001493"
$ws.Range("E28").Value = "'1.84%"
$ws.Range("D40").Value = "'0.03822"
$ws.Range("E40").Value = "'0.33%"
$ws.Range("D41").Value = "'0.006087"
$ws.Range("E41").Value = "'5.04%"
$ws.Range("E42").Value = "'0.23%"
$ws.Range("D43").Value = "'0.002299"
$ws.Range("E43").Value = "'-1.96%"
$ws.Range("D44").Value = "'0.01116"
$ws.Range("E44").Value = "'15.02%"
$ws.Range("D45").Value = "'0.00005494"
$ws.Range("E45").Value = "'4.87%"
$ws.Range("D46").Value = "'0.00000000749"
$ws.Range("E46").Value = "'-0.17%"
$ws.Range("D47").Value = "'0.08843"
$ws.Range("E47").Value = "'-4.94%"
$ws.Range("D48").Value = "'0.002414"
$ws.Range("E48").Value = "'12.55%"
$ws.Range("D49").Value = "'0.00002098"
$ws.Range("E49").Value = "'-0.17%"
$ws.Range("D50").Value = "'0.0001998"
$ws.Range("E50").Value = "'-0.17%"
